$d = $word.ActiveDocument

# Locate the "Version 1." text so we can compute character offsets robustly.
$find = $d.Content.Duplicate
$found = $find.Find.Execute("Version 1.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'Version 1.' text in document"
}
$base = $find.Start

# Character layout relative to $base:
#   0123456 7 8 9
#   Version   1  .
# "Version" = base+0 .. base+7
# " 1."     = base+7 .. base+10

# Step 1: split the run "Version" into "Versi" + "on" (matches target run
# boundary) by toggling a character-format property on the "Versi" sub-range
# and back, which forces Word to give that sub-range its own run without
# changing its visible formatting.
$rVersi = $d.Range($base, $base + 5)
$rVersi.Font.Bold = 1
$rVersi.Font.Bold = 0

# Step 2: change the version number "1" -> "2".
$rNum = $d.Range($base + 8, $base + 9)
$rNum.Text = "2"

# Step 3: remove the trailing "." (it currently sits before the _GoBack
# bookmark; we'll re-add a "." after the bookmark so the bookmark ends up
# sandwiched between " 2" and "." exactly as in the target).
$rDot = $d.Range($base + 9, $base + 10)
$rDot.Delete()

# Step 4: insert a new "." run right after the (now relocated) bookmark,
# which currently collapses to the end of "Version 2".
$rAfterBookmark = $d.Range($base + 9, $base + 9)
$rAfterBookmark.InsertAfter(".")
